$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-InCell($cell, $oldText, $newText) {
    $rng = $cell.Range
    $rng.Find.ClearFormatting()
    $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null
}

# --- Row 2 (Cluster 1, color D62728) ---
Replace-InCell $t.Cell(2,2) "emotion, behavior, adult, perception, affect, stress, threat, adolescent, depression, cardiovascular, social, reduce, heart rate variability, increase, startle, fear, anxiety" "perception, cognition, learning, affect, threat, stress, context, cardiovascular, social, pain, impact, feedback, fear, startle, cardiac, heart rate variability, interaction"
Replace-InCell $t.Cell(2,4) "2520 (1)" "2529 (1)"
Replace-InCell $t.Cell(2,5) "2034 (1)" "2088 (2)"

# --- Row 3 (Cluster 2, color 2CA02C) ---
Replace-InCell $t.Cell(3,2) "electrophysiological, inhibition, learning, reward, child, pain, performance, time, error, feedback, difference" "erp, eeg, attention, electrophysiological, auditory, p300, human, model, component, inhibition, performance, detection, index, control, individual"
Replace-InCell $t.Cell(3,3) "11" "15"
Replace-InCell $t.Cell(3,4) "2149 (2)" "2188 (2)"
Replace-InCell $t.Cell(3,5) "1850 (2)" "2182 (1)"

# --- Row 4 (Cluster 3, color 1F77B4) ---
Replace-InCell $t.Cell(4,2) "erp, attention, eeg, visual, auditory, p300, impact, human, component, model, detection" "emotion, neural, behavior, adult, memory, adolescent, depression, reduce, increase, reward, child, anxiety"
Replace-InCell $t.Cell(4,3) "11" "12"
Replace-InCell $t.Cell(4,4) "2077 (3)" "2139 (3)"
Replace-InCell $t.Cell(4,5) "1849 (3)" "1499 (4)"

# --- Row 5 (Cluster 4, color BCBD22) absorbs numeric values from the
#     about-to-be-deleted Row 6 (old Cluster 5, color 9467BD) ---
Replace-InCell $t.Cell(5,2) "neural, memory, brain, context, anticipation, mechanism, oscillations" "visual, brain, dynamic, oscillations, anticipation, cortex"
Replace-InCell $t.Cell(5,3) "7" "6"
Replace-InCell $t.Cell(5,4) "1642 (4)" "1251 (4)"
Replace-InCell $t.Cell(5,5) "1586 (4)" "1637 (3)"

# Row 5 becomes the new last row of the table, so its first two cells
# (which previously had no bottom border) now get the heavy bottom rule
# that used to sit on the (about to be removed) last row.
$bottom1 = $t.Cell(5,1).Borders.Item(-3)
$bottom1.LineStyle = 1
$bottom1.LineWidth = 6
$bottom1.Color = 6710886

$bottom2 = $t.Cell(5,2).Borders.Item(-3)
$bottom2.LineStyle = 1
$bottom2.LineWidth = 6
$bottom2.Color = 6710886

# --- Remove old Row 6 (Cluster 5, color 9467BD) entirely ---
$t.Rows.Item(6).Delete()
